# Update "想去人数" (F column) values across sheets to reflect the latest
# generated snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 1569
$ws.Range("F4").Value  = 819
$ws.Range("F5").Value  = 228
$ws.Range("F6").Value  = 62
$ws.Range("F7").Value  = 1127
$ws.Range("F8").Value  = 728
$ws.Range("F10").Value = 1398
$ws.Range("F11").Value = 277
$ws.Range("F12").Value = 1026
$ws.Range("F15").Value = 188
$ws.Range("F24").Value = 237

# --- Sheet: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 8
$ws.Range("F3").Value = 994
$ws.Range("F5").Value = 259
$ws.Range("F8").Value = 64

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 214

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 214
$ws.Range("F4").Value  = 1569
$ws.Range("F5").Value  = 8
$ws.Range("F6").Value  = 819
$ws.Range("F7").Value  = 228
$ws.Range("F8").Value  = 994
$ws.Range("F9").Value  = 62
$ws.Range("F10").Value = 1127
$ws.Range("F11").Value = 728
$ws.Range("F13").Value = 1398
$ws.Range("F14").Value = 277
$ws.Range("F15").Value = 1026
$ws.Range("F18").Value = 188
$ws.Range("F23").Value = 259
$ws.Range("F32").Value = 237
$ws.Range("F33").Value = 64

$wb.Save()
